# Agrupamento de distribuidoras: ENF+EMG=EMR e EBO+EPB=EPB
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count()

# Locate the rows for the distributors that need to be merged (search column A)
$rowEBO = 0
$rowEMG = 0
$rowENF = 0
$rowEPB = 0

for ($r = 2; $r -le $lastRow; $r++) {
    $name = $ws.Cells.Item($r, 1).Value()
    if ($name -eq "EBO") { $rowEBO = $r }
    if ($name -eq "EMG") { $rowEMG = $r }
    if ($name -eq "ENF") { $rowENF = $r }
    if ($name -eq "EPB") { $rowEPB = $r }
}

# Capture the underlying numeric values (casa_propria = col B, domicilios_pp = col C)
$eboB = $ws.Cells.Item($rowEBO, 2).Value()
$eboC = $ws.Cells.Item($rowEBO, 3).Value()
$epbB = $ws.Cells.Item($rowEPB, 2).Value()
$epbC = $ws.Cells.Item($rowEPB, 3).Value()
$enfB = $ws.Cells.Item($rowENF, 2).Value()
$enfC = $ws.Cells.Item($rowENF, 3).Value()
$emgB = $ws.Cells.Item($rowEMG, 2).Value()
$emgC = $ws.Cells.Item($rowEMG, 3).Value()

# New combined totals: EBO + EPB -> EPB ; ENF + EMG -> EMR
$newEpbB = $eboB + $epbB
$newEpbC = $eboC + $epbC
$newEpbD = $newEpbB / $newEpbC

$newEmrB = $enfB + $emgB
$newEmrC = $enfC + $emgC
$newEmrD = $newEmrB / $newEmrC

# Delete the four source rows, starting from the bottom-most row so that the
# remaining row indices captured above stay valid while we work our way up.
$rowsToDelete = @($rowEBO, $rowEMG, $rowENF, $rowEPB) | Sort-Object -Descending
foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete()
}

# The filtered/table range stays anchored to the rows that remain from the
# original table (i.e. it is not extended to cover the two newly appended
# rows below it). Capture the row count and (re)apply the AutoFilter now,
# BEFORE appending the merged rows below, so that the filter range does not
# grow to swallow the new rows underneath it.
$filterLastRow = $ws.UsedRange.Rows.Count()

$ws.AutoFilterMode = $false
$ws.Range("A1:D" + $filterLastRow).AutoFilter()

# Keep the workbook-level _FilterDatabase defined name in sync as well
for ($i = 1; $i -le $wb.Names.Count(); $i++) {
    $n = $wb.Names.Item($i)
    if ($n.Name() -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$D`$" + $filterLastRow
    }
}

# Now append the two merged rows at the bottom of the table
$emrRow = $filterLastRow + 1
$ws.Cells.Item($emrRow, 1).Value = "EMR"
$ws.Cells.Item($emrRow, 2).Value = $newEmrB
$ws.Cells.Item($emrRow, 3).Value = $newEmrC
$ws.Cells.Item($emrRow, 4).Value = $newEmrD

$epbRow = $filterLastRow + 2
$ws.Cells.Item($epbRow, 1).Value = "EPB"
$ws.Cells.Item($epbRow, 2).Value = $newEpbB
$ws.Cells.Item($epbRow, 3).Value = $newEpbC
$ws.Cells.Item($epbRow, 4).Value = $newEpbD

$ws.Range("A1").Select()
